$wb = $excel.ActiveWorkbook

# Map of worksheet name -> row => new value for column F ("想去人数")
$updates = @{
    2  = 1621
    3  = 9006
    6  = 691
    7  = 370
    8  = 184
    12 = 58
    15 = 4286
    20 = 332
    22 = 246
    24 = 2653
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
